$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.937.83'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.188.93'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.74'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.09'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +3.80%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("E8").Value = '  -3.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.32'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("E11").Value = '  -2.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.738.60'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.138'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -2.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.72'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000171'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.939.73'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +2.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.197.70'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +1.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.22'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.24'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.19'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '369.01'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -2.21%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  -1.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.58'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -0.57%  '

$ws.Range("E25").Value = '  +1.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.59'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +3.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0875'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.46'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.08'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.28'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +4.54%  '

$ws.Range("E34").Value = '  +2.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.90'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.26'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +5.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.785.46'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +5.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0708'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +2.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0309'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +6.97%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("E42").Value = '  -1.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.91'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +2.05%  '

$ws.Range("E44").Value = '  -0.78%  '

$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.229.30'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +0.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.980'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.15'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.57'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.794'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("E51").Value = '  +0.01%  '
